# Apply the "Study 1 almost done" edits to "Stimuli For Study.xlsx"
# Renames the generic "Choice 1" / " Choice 2" column headers (row 2, repeated
# across all five question blocks) to "Option 1" / " Option 2", and renames the
# Study 3 (Shipping) row labels "Date Received Days" / "Date Received Ratios"
# to "Shipping Time Days" / "Shipping Time Ratios".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 header cells: "Choice 1" -> "Option 1"
$ws.Range("B2").Value2 = "Option 1"
$ws.Range("E2").Value2 = "Option 1"
$ws.Range("H2").Value2 = "Option 1"
$ws.Range("K2").Value2 = "Option 1"
$ws.Range("N2").Value2 = "Option 1"

# Row 2 header cells: " Choice 2" -> " Option 2" (note the leading space)
$ws.Range("C2").Value2 = " Option 2"
$ws.Range("F2").Value2 = " Option 2"
$ws.Range("I2").Value2 = " Option 2"
$ws.Range("L2").Value2 = " Option 2"
$ws.Range("O2").Value2 = " Option 2"

# Study 3 (Shipping) row labels
$ws.Range("A15").Value2 = "    Shipping Time Days"
$ws.Range("A16").Value2 = "    Shipping Time Ratios"

# Update the saved selection to match the author's workbook state
$ws.Range("Q9").Select()
